$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 32 (shifts rows 32.. down by one, preserving their
# content/formatting exactly, and gives the new row inherited formatting).
$ws.Rows("32:32").Insert()

# Populate the newly inserted row 32 with the new localization entry.
$ws.Cells.Item(32, 2).Value = 'localization\strings'
$ws.Cells.Item(32, 3).Value = 'strWindowPos'
$ws.Cells.Item(32, 4).Value = 'In "settings" form, tab "User interface"'
$ws.Cells.Item(32, 5).Value = 'Remember window position and size on startup'

# The existing "strChkDlgPath" row (row 25) also gets the same comment applied
# retroactively (it shares the same settings-form/"User interface" tab context).
$ws.Cells.Item(25, 4).Value = 'In "settings" form, tab "User interface"'

# Grow the table (ListObject) by one row to match the new data extent.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:F204"))

# Widen column D slightly to accommodate the new, longer comment text.
$ws.Columns("D").ColumnWidth = 34.83
